$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 data (Day 17) ---
# (Order matches the original author's edit so new shared-string entries
# come out in the same sequence: Day 17, then the three filenames, then
# the topics cell.)
$ws.Range("A18").Value = "Day 17"
$ws.Range("B18").Value = 45819

# --- Hyperlinks for the new "Day 17" problems ---
# (TextToDisplay is passed as the full URL so the saved "display" attribute
# matches it; the cell's Value is then set explicitly afterwards so the
# visible cell text stays as the short filename.)
$ws.Hyperlinks.Add($ws.Range("C18"), "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 17/Minimum Size Subarray Sum.py", "", "Minimum Size Subarray Sum.py", "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 17/Minimum Size Subarray Sum.py") | Out-Null
$ws.Range("C18").Value = "Minimum Size Subarray Sum.py"

$ws.Hyperlinks.Add($ws.Range("D18"), "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 17/Permutation in String.py", "", "Permutation in String.py", "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 17/Permutation in String.py") | Out-Null
$ws.Range("D18").Value = "Permutation in String.py"

$ws.Hyperlinks.Add($ws.Range("E18"), "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 17/Search Insert Position.py", "", "Search Insert Position.py", "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 17/Search Insert Position.py") | Out-Null
$ws.Range("E18").Value = "Search Insert Position.py"

$ws.Range("F18").Value = "Binary Search, Sliding window"
$ws.Range("G18").Value = "S"

# --- Re-apply the same formatting used by the rows above (row 17) ---
# Date style for B18 (built-in date number format, same xf as B17)
$ws.Range("B17").Copy()
$ws.Range("B18").PasteSpecial(-4122)

# Hyperlink cell style for C18/D18/E18 (same xf as C17/D17/E17), applied
# *after* Hyperlinks.Add so it isn't overwritten by the default link style
$ws.Range("C17").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E17").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- Update the active selection shown in the sheet view ---
$ws.Range("D12").Select()
